$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new value is a plain decimal must be pre-formatted as Text
# so Excel keeps the exact original digit grouping/trailing zeros instead of
# silently coercing the literal into a Double (source data are inline strings,
# not numbers). The COM shim only honours NumberFormat on a single contiguous
# area, so each block below is set individually.
$ws.Range("D4:D5").NumberFormat = "@"
$ws.Range("D7:D8").NumberFormat = "@"
$ws.Range("D11:D16").NumberFormat = "@"
$ws.Range("D18:D20").NumberFormat = "@"
$ws.Range("D22:D27").NumberFormat = "@"
$ws.Range("D29:D51").NumberFormat = "@"

# Rows 2-33: refresh whichever Price (D) / Volume 1h (E) cells actually moved
$ws.Range("D2").Value = "26.727.30"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.727.17"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").Value = "0.9973"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "242.19"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "0.4929"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.2619"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "1.730.21"
$ws.Range("D11").Value = "15.91"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "0.06989"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "0.6124"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "4.515"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "77.24"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "0.9978"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "26.523.03"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "0.9970"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "0.000007204"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "1.949.02"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "4.454"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "8.590"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "5.113"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "138.24"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "15.34"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "1.758"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "106.32"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "3.918"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").Value = "0.07989"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "3.672"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "0.04504"
$ws.Range("E33").Value = "  -0.74%  "

# Rows 34-51: "Frax" dropped from the ranking, every following coin shifted up
# one row, and "NEARProtocol" newly appended at row 51.
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.607"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.004"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6279"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.9339"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.032"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.416"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "0.9988"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01517"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.603"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "99.60"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3868"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "6.912"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1159"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05381"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "7.851"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.35"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "51.81"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.236"
$ws.Range("E51").Value = "  -0.90%  "
